# Edit the "Some notes, for this document I pushed..." paragraph:
#   1. change "document" -> "blog" (splitting that run into 3 pieces)
#   2. color every run in that paragraph (and the paragraph mark) red (FF0000)

$d = $word.ActiveDocument

# --- Step 1: text substitution --------------------------------------------
# Locate + replace "document" (scoped by surrounding words so we hit the one
# and only occurrence inside this paragraph) with "blog".
$tr = $d.Content
$foundText = $tr.Find.Execute("document I pushed our IoT data stream into", $true, $false, $false, $false, $false, $true, 1, $false, "blog I pushed our IoT data stream into", 2)
if (-not $foundText) {
    throw "Could not find target sentence to update"
}

# --- Step 2: locate the (now merged-into-one-run) paragraph via its prefix
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("Some notes, for this blog I pushed our IoT data stream into ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundAnchor) {
    throw "Could not find anchor run after text substitution"
}
$paraStart = $anchor.Start
$targetParagraph = $anchor.Paragraphs(1)

# --- Step 3: apply red color segment-by-segment, recreating the run breaks
$segments = @(
    "Some notes, for this ",
    "blog",
    " I pushed our IoT data stream into ",
    "Apache ",
    "Kafka",
    " topics",
    ", going forward Fluss will ",
    "be able to ",
    "present a ",
    "Apache ",
    "Kafka compatible endpoint allowing for data to be published directly into Fluss tables ",
    "(using Kafka protocol) ",
    "which will simplify our stack and result in much " + [string][char]0x201C + "fresher" + [string][char]0x201D + " data for analytics",
    " and lower cost as we will have less technology involved."
)

$pos = $paraStart
foreach ($seg in $segments) {
    $segLen = $seg.Length
    $rng = $d.Range($pos, $pos + $segLen)
    $rng.Font.Color = 255
    $pos = $pos + $segLen
}

# --- Step 4: color the paragraph mark too (pPr/rPr) -----------------------
# Re-applying color across the whole paragraph (content + mark) keeps the
# already-split runs intact (they already carry the colour) while adding the
# colour to the paragraph-mark run properties that sit in w:pPr/w:rPr.
$targetParagraph.Range.Font.Color = 255

Write-Host "Done. foundText=$foundText foundAnchor=$foundAnchor paraStart=$paraStart finalPos=$pos"
